# backwardElimination.xlsx : refresh the "Summary" regression printouts on
# each of the 3 sheets so the embedded statsmodels OLS report shows the
# latest run's Date/Time stamp (Wed 01 Jan 2020 23:18:49 -> Thu 02 Jan 2020
# 20:48:41). Everything else in the report text is left untouched.

$wb = $excel.ActiveWorkbook

$oldDate = "Date:                Wed, 01 Jan 2020"
$newDate = "Date:                Thu, 02 Jan 2020"
$oldTime = "Time:                        23:18:49"
$newTime = "Time:                        20:48:41"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = $cell.Value2

    if ($text -and $text.Contains($oldDate)) {
        $text = $text.Replace($oldDate, $newDate)
        $text = $text.Replace($oldTime, $newTime)
        $cell.Value = $text
    }
}
